$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "61.392.77"
$ws.Cells.Item(2, 5).Value = "  -5.99%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.976.45"
$ws.Cells.Item(3, 5).Value = "  -7.38%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.01%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'541.65"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -5.98%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'152.85"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -8.71%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.01%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.561"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -6.25%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.978.53"
$ws.Cells.Item(9, 5).Value = "  -7.03%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -7.28%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -8.12%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -6.99%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.502.13"
$ws.Cells.Item(13, 5).Value = "  -7.10%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -3.82%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "61.571.07"
$ws.Cells.Item(15, 5).Value = "  -5.71%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'23.58"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -8.20%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.985.74"
$ws.Cells.Item(17, 5).Value = "  -6.97%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -7.33%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'388.47"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -6.16%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -4.87%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'11.86"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -8.26%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.60"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -8.09%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.09%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'64.69"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -7.25%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -5.03%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.186"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -8.57%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.11%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "0.0₃0934"
$ws.Cells.Item(28, 5).Value = "  -11.29%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'8.37"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -6.28%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'1.00"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.00%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -7.34%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'20.25"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -6.55%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'159.23"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.23%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'5.99"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -6.59%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -8.13%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -7.18%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -6.97%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -9.70%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "2.432.77"
$ws.Cells.Item(39, 5).Value = "  -11.60%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Filecoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40, 4).Value = "'3.87"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -6.84%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "OKB"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(41, 4).Value = "'37.06"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -5.32%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -8.44%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -8.28%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -6.69%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.999"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.09%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.0245"
$ws.Cells.Item(46, 4).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 4).Value = "'4.89"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -13.05%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0954"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.80%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(49, 4).Value = "'19.63"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -9.14%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(50, 4).Value = "'10.47"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.00%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'263.17"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -11.48%  "
